$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-10-26 Sunday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-10-27 Monday", 2) | Out-Null

# Update the division problems in the worksheet table, addressed by
# (row, column) so that duplicate problem text (e.g. "84÷4=", "90÷4=")
# is replaced unambiguously and in document order.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "79÷6="
$t.Cell(1, 2).Range.Text  = "97÷6="
$t.Cell(1, 3).Range.Text  = "77÷2="
$t.Cell(1, 4).Range.Text  = "47÷9="
$t.Cell(1, 5).Range.Text  = "43÷6="

$t.Cell(5, 1).Range.Text  = "44÷2="
$t.Cell(5, 2).Range.Text  = "38÷4="
$t.Cell(5, 3).Range.Text  = "28÷5="
$t.Cell(5, 4).Range.Text  = "99÷4="
$t.Cell(5, 5).Range.Text  = "37÷6="

$t.Cell(9, 1).Range.Text  = "35÷3="
$t.Cell(9, 2).Range.Text  = "63÷5="
$t.Cell(9, 3).Range.Text  = "54÷8="
$t.Cell(9, 4).Range.Text  = "25÷5="
$t.Cell(9, 5).Range.Text  = "25÷6="

$t.Cell(13, 1).Range.Text = "72÷9="
$t.Cell(13, 2).Range.Text = "21÷7="
$t.Cell(13, 3).Range.Text = "13÷3="
$t.Cell(13, 4).Range.Text = "90÷7="
$t.Cell(13, 5).Range.Text = "31÷3="

$t.Cell(17, 1).Range.Text = "43÷8="
$t.Cell(17, 2).Range.Text = "15÷6="
$t.Cell(17, 3).Range.Text = "89÷7="
$t.Cell(17, 4).Range.Text = "44÷5="
$t.Cell(17, 5).Range.Text = "25÷8="
